# Quarterly income-statement refresh for سهگمت-سیمان‌هگمتان‌ (Overview sheet)
# Drops the two oldest quarters (columns D:E) and appends two new quarters
# (columns L:M), matching how the source workbook was rolled forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- 1. Drop the two oldest quarter columns, shifting everything left ----
$ws.Range("D1:E28").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# --- 2. Seed the two freed-up trailing columns (L:M) with K's formatting -
$ws.Range("J1:K28").Copy()
$ws.Range("L1:M28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Match the custom column widths that shift along with the data (31-char
# wide "first quarter of a bank" columns land on D, H, L; the rest stay 29).
$ws.Range("L1").ColumnWidth = 30.1666666666667
$ws.Range("M1").ColumnWidth = 28.1666666666667

# --- 3. New quarter headers (row 8) and publish dates (row 9) ------------
$ws.Range("L8").Value2 = "فصل چهارم منتهی به 1401/10"
$ws.Range("M8").Value2 = "فصل اول منتهی به 1402/01"

$ws.Range("H9").Value2 = "1402-02-13 (9)"
$ws.Range("I9").Value2 = "1402-02-29 (2)"
$ws.Range("L9").Value2 = "1402-02-29 (3)"
$ws.Range("M9").Value2 = "1402-02-29"

# --- 4. Refresh the financial figures for the two new columns, and fix up
#        a handful of historical values that were recomputed under the new
#        "read_price" algorithm (commit message) -------------------------
$ws.Range("L11").Value2 = 1962122
$ws.Range("M11").Value2 = 3019253

$ws.Range("L12").Value2 = -1458189
$ws.Range("M12").Value2 = -1284053

$ws.Range("L13").Value2 = 503933
$ws.Range("M13").Value2 = 1735200

$ws.Range("L14").Value2 = -156373
$ws.Range("M14").Value2 = -225800

$ws.Range("L15").Value2 = 0
$ws.Range("M15").Value2 = 0

$ws.Range("L16").Value2 = 4253
$ws.Range("M16").Value2 = 19526

$ws.Range("L17").Value2 = 351813
$ws.Range("M17").Value2 = 1528926

$ws.Range("L18").Value2 = -1886
$ws.Range("M18").Value2 = 0

$ws.Range("I19").Value2 = 0
$ws.Range("L19").Value2 = 757483
$ws.Range("M19").Value2 = 102205

$ws.Range("I20").Value2 = 678804
$ws.Range("L20").Value2 = 1107410
$ws.Range("M20").Value2 = 1631131

$ws.Range("L21").Value2 = 115923
$ws.Range("M21").Value2 = -164216

$ws.Range("I22").Value2 = 570609
$ws.Range("L22").Value2 = 1223333
$ws.Range("M22").Value2 = 1466915

$ws.Range("L23").Value2 = 0
$ws.Range("M23").Value2 = 0

$ws.Range("I24").Value2 = 570609
$ws.Range("L24").Value2 = 1223333
$ws.Range("M24").Value2 = 1466915

$ws.Range("I25").Value2 = 783
$ws.Range("L25").Value2 = 1679
$ws.Range("M25").Value2 = 2013

$ws.Range("L26").Value2 = 728789
$ws.Range("M26").Value2 = 728789

$ws.Range("I27").Value2 = 783
$ws.Range("L27").Value2 = 1679
$ws.Range("M27").Value2 = 2013
